# Apply "Test data Updation" edits described by the commit diff.
#
# Summary of changes:
#  - "Capabilities" sheet: EndPoint URL updated from the old pCloudy device
#    endpoint to the new US endpoint (rows 2 & 3, column D).
#  - "DeviceList" sheet: device/version/OS rows swapped from an Android
#    device (OnePlus) to iOS devices (iPhone 11 Pro / iPhone 8), and the
#    OperatingSystem row switched from pCloudyAndroid to pCloudyIOS.
#  - Selection/active-cell bookmarks updated on both sheets to match the
#    state the workbook was left in when saved.

$wb = $excel.ActiveWorkbook

$capsSheet = $wb.Worksheets.Item("Capabilities")
$deviceSheet = $wb.Worksheets.Item("DeviceList")

# ---- "DeviceList" sheet: first (B) device/version columns ----------------
# Device identifier -> Apple iPhone 11 Pro (was a OnePlus Android device).
$deviceSheet.Range("B1").Value = "APPLE_iPhone11Pro_iOS_14.4.0_6ccce"
# Version number matching the new device.
$deviceSheet.Range("B2").Value = "14.4.0"

# ---- "Capabilities" sheet --------------------------------------------------
# EndPoint column (D) for the two credential rows now points at the US
# pCloudy endpoint instead of the generic "device" one.
$capsSheet.Range("D2").Value = "https://us.pcloudy.com"
$capsSheet.Range("D3").Value = "https://us.pcloudy.com"

# ---- "DeviceList" sheet: second (C) device/version columns ---------------
# Device identifier -> Apple iPhone 8 (was a OnePlus Android device).
$deviceSheet.Range("C1").Value = "APPLE_iPhone8_iOS_14.1.0_81551"
# Version number matching the new device.
$deviceSheet.Range("C2").Value = "14.1.0"

# Row 3: OperatingSystem switches from Android to iOS pCloudy platform.
$deviceSheet.Range("B3").Value = "pCloudyIOS"
$deviceSheet.Range("C3").Value = "pCloudyIOS"

# Leave the selection where the author left it when the file was saved.
$capsSheet.Range("D7").Select()
$deviceSheet.Range("C18").Select()
